# Chop.Calc.xlsx - rename "Index" column to "i" and renumber it 0-based.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell (this also renames the table column "Index" -> "i",
# and causes "Index" to drop out of the shared string table while "i" is
# appended as a new shared string).
$ws.Range("A1").Value = "i"

# Renumber column A from 1-based (1..502) to 0-based (0..501) for all data rows.
for ($r = 2; $r -le 503; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}

# Narrow column A now that it holds shorter values ("i" header + up to 3 digits).
$col = $ws.Columns.Item(1)
$col.ColumnWidth = 3.17
